# Apply updated crypto price / 1h-volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.797.20'
$ws.Range("E2").Value = '  +1.20%  '
$ws.Range("D3").Value = '1.888.10'
$ws.Range("E3").Value = '  +1.84%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'239.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = "'0.4759"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.68%  '
$ws.Range("D8").Value = "'0.2884"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.91%  '
$ws.Range("D9").Value = "'0.06595"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.33%  '
$ws.Range("E10").Value = '  +10.22%  '
$ws.Range("D11").Value = "'99.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +17.73%  '
$ws.Range("D12").Value = '1.886.77'
$ws.Range("E12").Value = '  +1.72%  '
$ws.Range("D13").Value = "'0.07613"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("E14").Value = '  +3.92%  '
$ws.Range("D15").Value = "'0.6634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.97%  '
$ws.Range("D16").Value = "'308.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +34.77%  '
$ws.Range("D17").Value = '30.797.76'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("E18").Value = '  +5.57%  '
$ws.Range("D19").Value = "'0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = "'0.000007598"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("D21").Value = '2.124.28'
$ws.Range("E21").Value = '  +1.70%  '
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = "'5.124"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'6.222"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.14%  '
$ws.Range("D25").Value = "'9.313"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").Value = "'167.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("D27").Value = "'20.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.09%  '
$ws.Range("D28").Value = "'1.951"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.68%  '
$ws.Range("E29").Value = '  +5.57%  '
$ws.Range("D30").Value = "'1.353"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.56%  '
$ws.Range("E31").Value = '  +2.14%  '
$ws.Range("D32").Value = "'3.988"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.19%  '
$ws.Range("D33").Value = "'0.05055"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.50%  '
$ws.Range("D34").Value = "'1.174"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.86%  '
$ws.Range("D35").Value = "'0.7301"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").Value = "'0.01959"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("D38").Value = "'2.703"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").Value = "'2.074"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.38%  '
$ws.Range("D40").Value = "'0.9071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.20%  '
$ws.Range("D41").Value = "'108.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("D42").Value = "'0.9999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = "'0.4217"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.00%  '
$ws.Range("D44").Value = "'5.654"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.44%  '
$ws.Range("D45").Value = "'7.408"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("E46").Value = '  +7.48%  '
$ws.Range("D47").Value = "'9.049"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.43%  '
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("E49").Value = '  +4.02%  '
$ws.Range("D50").Value = "'0.05635"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.51%  '
$ws.Range("E51").Value = '  +2.18%  '
